$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("work")

# Apply the same per-column cell style used by the existing data rows (row 23, a
# full A-F row) to the new rows BEFORE setting values, by copying single cells'
# formats - this avoids Excel re-deriving/duplicating number-format codes.
$ws.Range("A23").Copy()
$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Copy()
$ws.Range("B27:B29").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Copy()
$ws.Range("C27:C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy()
$ws.Range("D27:D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Copy()
$ws.Range("E27:E28").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Copy()
$ws.Range("F27:F28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 27
$ws.Range("A27").Value = 43913
$ws.Range("B27").Value = "1"
$ws.Range("C27").Value = "Calculate sphere distance between two points, calculate distance between a point to a segment"
$ws.Range("D27").Value = 0.70138888888888884
$ws.Range("E27").Value = 0.74305555555555547
$ws.Range("F27").Value = 1

# Row 28
$ws.Range("B28").Value = "2"
$ws.Range("C28").Value = "Find nearest Link to accident spot - 2: Format path.json"
$ws.Range("D28").Value = 0.75
$ws.Range("E28").Value = 0.91666666666666663
$ws.Range("F28").Value = 4

# Row 29
$ws.Range("B29").Value = "3"

# Update selection to match final state
$ws.Range("C29").Select()
